$wb = $excel.ActiveWorkbook

# Rename Sheet2 -> loginDetails
$wb.Worksheets.Item("Sheet2").Name = "loginDetails"

$ws1 = $wb.Worksheets.Item("createAccount")

# Remove the hyperlink + formatting on D2, then set its new value
$ws1.Range("D2").Hyperlinks.Delete()
$ws1.Range("D2").ClearFormats()
$ws1.Range("D2").Value = "johnsmith10@live.com"

# Apply Text number format to the range that was reformatted, and give U2 a
# text value instead of a number
$ws1.Range("F2").NumberFormat = "@"
$ws1.Range("H2").NumberFormat = "@"
$ws1.Range("O2").NumberFormat = "@"
$ws1.Range("R2").NumberFormat = "@"
$ws1.Range("S2").NumberFormat = "@"
$ws1.Range("T2").NumberFormat = "@"
$ws1.Range("U2").NumberFormat = "@"
$ws1.Range("U2").Value = "#55592111880"

# Update selection on loginDetails, then switch back so createAccount stays
# the active/selected sheet (as it was before the edit)
$ws2 = $wb.Worksheets.Item("loginDetails")
$ws2.Range("J28").Select()

$ws1.Range("D7").Select()
